$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Property1" to "DataNode" (part of unifying the
# DataNode / DataTable / Entity naming scheme referenced in the commit).
$ws.Name = "DataNode"

# Re-select the previously selected cell in the frozen (bottom-left) pane,
# moving the selection from D9 to F25.
$ws.Range("F25").Select()
